# Apply the text corrections on the "Business Approved List" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Approved List")

# Fix typo: "medicate" -> "medicare"
$ws.Range("B4").Value = "Group. medicare (V01)_- ESA (PPO)-Local 147 -Construction-WorkersFund"

# Drop trailing "_-Medicare" suffix
$ws.Range("B9").Value = "Group.Medicare (SO3) (HMO) Distict-Council33"

# Update the active selection shown when the sheet is opened
$ws.Activate() | Out-Null
$ws.Range("B15").Select() | Out-Null
